# Update "想去人数" (F column) values on the "展览" and "全部类型" sheets.
# Row 3: 50 -> 52
# Row 4: 256 -> 258
# Row 5: 3971 -> 3983

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F3").Value = 52
    $ws.Range("F4").Value = 258
    $ws.Range("F5").Value = 3983
}
